# "Đơn phụ phẫu 2" is the 6th tab in this workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Đơn phụ phẫu 2")

# Columns that hold numeric totals (Đơn giá gốc, Upsale, Đơn giá, Thanh toán
# lần đầu, Trả sau, Đã thanh toán, Dư nợ) plus the service-code column (Mã
# dịch vụ) become 0 on the new blank data row.
$zeroCols = @("B", "I", "K", "L", "M", "N", "O", "P")
foreach ($col in $zeroCols) {
    $ws.Range($col + "2").Value = 0
}

# Remaining (text) columns stay blank — an empty string assignment would be
# treated as "clear the cell" (same as real Excel), so an empty formula is
# used to materialize an actual blank cell on row 2.
$textCols = @("A", "C", "D", "E", "F", "G", "H", "J", "Q", "R", "S", "T")
foreach ($col in $textCols) {
    $ws.Range($col + "2").Formula = "=""""" 
}
